# KPI_APPLE.xlsx update:
#  - Header row renamed: A1 "KPI_Microsoft" -> "Fecha", B1 "SUMA" -> "Ventas"
#  - New column C added with header "kpi"
#  - C2 gets a formula that blends the current and prior year's ROE,
#    formatted as a percentage (0.00%)
#  - Selection left on C3, matching the saved cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Ventas"

# Add the new "kpi" column
$ws.Range("C1").Value = "kpi"
$ws.Range("C2").Formula = "=B2+B3/2"
$ws.Range("C2").NumberFormat = "0.00%"

# Match the selection recorded in the saved workbook
$ws.Range("C3").Select() | Out-Null
